# "changed textfeilds to buttons for somepages"
# Restructure several sheets: drop stray extra data rows, blank out some
# row-2 values, and split several header rows into per-"term" pairs
# (each original header "X" becomes "X term 1" and "X term 2").

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "Sheet": currently a single blank cell (A1). Extend the used
# range down to A2 (still blank) without changing its value.
# ---------------------------------------------------------------------
$wsMain = $wb.Worksheets.Item("Sheet")
$wsMain.Range("A1").Font.Bold = $false
$wsMain.Range("A2").Font.Bold = $false

# ---------------------------------------------------------------------
# "cover_page": drop the "Varad" row (row 3) entirely, and blank out
# the "Manoj" row (row 2) while keeping the Admission Number (12).
# ---------------------------------------------------------------------
$wsCover = $wb.Worksheets.Item("cover_page")
$wsCover.Rows.Item(3).Delete()
$wsCover.Range("A2").Value = ""
$wsCover.Range("C2:J2").Value = ""
$wsCover.Range("A2:J2").Font.Bold = $false
$wsCover.Range("A1").Select()

# ---------------------------------------------------------------------
# "first_page": drop the "Manoj" row (row 3) entirely, and leave row 2
# present but blank across all of its columns.
# ---------------------------------------------------------------------
$wsFirst = $wb.Worksheets.Item("first_page")
$wsFirst.Rows.Item(3).Delete()
$wsFirst.Range("A2:J2").Value = ""
$wsFirst.Range("A2:J2").Font.Bold = $false
$wsFirst.Range("A1").Select()

# ---------------------------------------------------------------------
# "HEALTH & WELLBEING (HW)": split each header into "term 1"/"term 2"
# pairs (10 headers -> 20 columns), clear the old sample row (row 2)
# and drop the extra blank row (row 3).
# ---------------------------------------------------------------------
$wsHW = $wb.Worksheets.Item("HEALTH & WELLBEING (HW)")
$hwHeaders = @("HW 4.10 term 1", "HW 4.10 term 2", "HW 4.11 term 1", "HW 4.11 term 2", "HW 4.7 term 1", "HW 4.7 term 2", "HW 5.8 term 1", "HW 5.8 term 2", "HW 5.13A term 1", "HW 5.13A term 2", "HW 5.13a term 1", "HW 5.13a term 2", "HW 5.16 term 1", "HW 5.16 term 2", "HW 5.17 term 1", "HW 5.17 term 2", "HW 5.18a term 1", "HW 5.18a term 2", "HW 5.18b term 1", "HW 5.18b term 2")
$wsHW.Rows.Item(3).Delete()
$wsHW.Range("A2").Value = ""
for ($i = 0; $i -lt $hwHeaders.Length; $i++) {
    $wsHW.Cells.Item(1, $i + 1).Value = $hwHeaders[$i]
}
$wsHW.Range("A2").Font.Bold = $false
$wsHW.Range("A1").Select()

# ---------------------------------------------------------------------
# "Effective Communication (ECL)": split each header into "term
# 1"/"term 2" pairs (11 headers -> 22 columns) and drop the extra
# blank row (row 3); row 2 stays present but blank.
# ---------------------------------------------------------------------
$wsECL = $wb.Worksheets.Item("Effective Communication (ECL)")
$eclHeaders = @("ECL 15.4 term 1", "ECL 15.4 term 2", "ECL 15.5a term 1", "ECL 15.5a term 2", "ECL 15.5b term 1", "ECL 15.5b term 2", "ECL 15.5c term 1", "ECL 15.5c term 2", "ECL 15.7 term 1", "ECL 15.7 term 2", "ECL 15.8 term 1", "ECL 15.8 term 2", "ECL 15.9 term 1", "ECL 15.9 term 2", "ECL 2-5.10 term 1", "ECL 2-5.10 term 2", "ECL 2-5.12 term 1", "ECL 2-5.12 term 2", "ECL 2-5.1a term 1", "ECL 2-5.1a term 2", "ECL 2-5.2 term 1", "ECL 2-5.2 term 2")
$wsECL.Rows.Item(3).Delete()
for ($i = 0; $i -lt $eclHeaders.Length; $i++) {
    $wsECL.Cells.Item(1, $i + 1).Value = $eclHeaders[$i]
}
$wsECL.Range("A2").Font.Bold = $false
$wsECL.Range("A1").Select()

# ---------------------------------------------------------------------
# "Involved Learners (IL)": split each header into "term 1"/"term 2"
# pairs (13 headers -> 26 columns) and drop the extra blank row
# (row 3); row 2 stays present but blank.
# ---------------------------------------------------------------------
$wsIL = $wb.Worksheets.Item("Involved Learners (IL)")
$ilHeaders = @("IL 4.1 term 1", "IL 4.1 term 2", "IL 4.2a term 1", "IL 4.2a term 2", "IL 4.8c term 1", "IL 4.8c term 2", "IL 4.6 term 1", "IL 4.6 term 2", "IL 4.11 term 1", "IL 4.11 term 2", "IL 4.13 term 1", "IL 4.13 term 2", "IL 4.16 term 1", "IL 4.16 term 2", "IL 4.9 term 1", "IL 4.9 term 2", "IL 4.20 term 1", "IL 4.20 term 2", "IL 4.25 term 1", "IL 4.25 term 2", "IL 4.27 term 1", "IL 4.27 term 2", "IL 4.29 term 1", "IL 4.29 term 2", "IL 4.30 term 1", "IL 4.30 term 2")
$wsIL.Rows.Item(3).Delete()
for ($i = 0; $i -lt $ilHeaders.Length; $i++) {
    $wsIL.Cells.Item(1, $i + 1).Value = $ilHeaders[$i]
}
$wsIL.Range("A2").Font.Bold = $false
$wsIL.Range("A1").Select()
